$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "F2" "9840049305"
Set-TextValue "F3" "9840040243"
Set-TextValue "F4" "9840009946"
Set-TextValue "F5" "9840041029"
Set-TextValue "F6" "9840067331"
Set-TextValue "F7" "9840037015"
Set-TextValue "F8" "9840029760"
Set-TextValue "F9" "9840002606"
Set-TextValue "F10" "9840060400"

Set-TextValue "AM2" "0"
Set-TextValue "AN2" "0"
Set-TextValue "AO2" "3"
